$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add two new rows to the table (extends range/dimension A1:E30 -> A1:E32)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Row 31: 647. Palindromic Substring ---
$ws.Range("A31").Value = "647. Palindromic Substring"
$ws.Range("E31").Value = "https://leetcode.com/problems/palindromic-substrings/solutions/105689/java-solution-8-lines-extendpalindrome/comments/1017849 "
$ws.Hyperlinks.Add($ws.Range("E31"), "https://leetcode.com/problems/palindromic-substrings/solutions/105689/java-solution-8-lines-extendpalindrome/comments/1017849") | Out-Null
$ws.Range("D31").Value = "2 passes - Odd and Even. We consider each index as mid and expand outwards, but also consider the next adjacent to capture the even ones. This is the general formula for finding palindromes. The Dynamic Programming solution uses the dp array to track inner windows."
$ws.Range("B31").Value = "Medium"
$ws.Range("C31").Value = "Dynamic Programming"

# --- Row 32: 36. Valid Sudoku ---
$ws.Range("A32").Value = "36. Valid Sudoku"
$ws.Range("D32").Value = "Use a Hash Set for each rule. For the subsquare rule, we use an array [i][j] to get the the subsquare. Divide by 3 on the row and column to get the coordinates and obtain the subsquare it is in."
$ws.Range("E32").Value = "https://leetcode.com/problems/valid-sudoku/solutions/15472/short-simple-java-using-strings/ "
$ws.Hyperlinks.Add($ws.Range("E32"), "https://leetcode.com/problems/valid-sudoku/solutions/15472/short-simple-java-using-strings/") | Out-Null
$ws.Range("B32").Value = "Medium"
$ws.Range("C32").Value = "Arrays"

# Copy cell formatting (fill on Difficulty column, Hyperlink style on Link column, etc.)
# from the previous last data row (30) onto the two newly added rows, without
# disturbing the values/hyperlinks already set above.
$ws.Range("A30:E30").Copy() | Out-Null
$ws.Range("A31:E31").PasteSpecial(-4122) | Out-Null
$ws.Range("A30:E30").Copy() | Out-Null
$ws.Range("A32:E32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Scroll the view down a bit (best effort, matches topLeftCell B7 -> B10)
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 2
